$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

$ws.Range("A2").Value = "Gato"
$ws.Range("A3").Value = "Arquero"
$ws.Range("A4").Value = "Ficha"
$ws.Range("A5").Value = "Poker"

$ws.Range("A5").Select()
